# Update the "dSF" column (F) values on Sheet1 to reflect the repulled /
# recalculated data, per the commit "repull data, push all data, mean
# calculation".
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 2
$ws.Range("F14").Value = 2
$ws.Range("F16").Value = -5
$ws.Range("F21").Value = -3
$ws.Range("F24").Value = 0
$ws.Range("F26").Value = -3
$ws.Range("F27").Value = 4
$ws.Range("F29").Value = 1
$ws.Range("F32").Value = -2
$ws.Range("F35").Value = -12
$ws.Range("F37").Value = -7
$ws.Range("F39").Value = -10
$ws.Range("F40").Value = -4
$ws.Range("F41").Value = -2
$ws.Range("F46").Value = -2

$wb.Save()
